$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2024-03-05 Tuesday"; New = "2024-03-06 Wednesday" },
    @{ Old = "476×2="; New = "682×6=" },
    @{ Old = "972×6="; New = "763×5=" },
    @{ Old = "606×8="; New = "158×7=" },
    @{ Old = "333×6="; New = "878×2=" },
    @{ Old = "211×8="; New = "782×9=" },
    @{ Old = "950×3="; New = "573×4=" },
    @{ Old = "585×5="; New = "169×4=" },
    @{ Old = "911×4="; New = "707×2=" },
    @{ Old = "785×9="; New = "979×6=" },
    @{ Old = "763×7="; New = "258×3=" },
    @{ Old = "292×2="; New = "933×5=" },
    @{ Old = "416×8="; New = "625×7=" },
    @{ Old = "702×3="; New = "490×2=" },
    @{ Old = "649×8="; New = "703×3=" },
    @{ Old = "240×7="; New = "577×2=" },
    @{ Old = "332×5="; New = "259×8=" },
    @{ Old = "602×8="; New = "860×9=" },
    @{ Old = "866×2="; New = "905×4=" },
    @{ Old = "111×2="; New = "679×8=" },
    @{ Old = "165×5="; New = "856×4=" },
    @{ Old = "427×8="; New = "646×4=" },
    @{ Old = "971×7="; New = "818×5=" },
    @{ Old = "909×3="; New = "713×3=" },
    @{ Old = "383×4="; New = "347×5=" },
    @{ Old = "889×4="; New = "457×3=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2)
}
